$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing API test-case rows (4-10); only the two
# newly-added front-end TCs remain below the header row.
$ws.Range("A4:C10").EntireRow.Delete()

# Row 2 -> iExamRegression front-end TC: Ranking Question
$ws.Range("A2").Value = "iEX_TC_ID_104"
$ws.Range("B2").Value = "@iExamRegression Validation of Exam Section > Ranking Question"
$ws.Range("C2").Value = "passed"

# Row 3 -> iExamRegression front-end TC: ISAWE CASE Question
$ws.Range("A3").Value = "iEX_TC_ID_107"
$ws.Range("B3").Value = "@iExamRegression Validation of Exam Section > ISAWE CASE Question"
$ws.Range("C3").Value = "passed"
